$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a cell's value as plain text, preserving the worksheet's
# default (unstyled) cell formatting. Many of the values here look like
# numbers (e.g. "0.7133", "1.001") but must be stored as text, matching
# the original workbook's inlineStr cells. Forcing NumberFormat to "@"
# (Text) before assignment prevents Excel from re-interpreting the
# string as a number/date, and resetting Style back to "Normal"
# afterwards avoids leaving a stray number-format style on the cell.
function Set-TextValue($Row, $Col, $Value) {
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.Style = "Normal"
}

$changes = @(
    @{ Row = 2; Col = 4; Letter = "D"; Value = "29.380.34" },
    @{ Row = 2; Col = 5; Letter = "E"; Value = "  +0.20%  " },
    @{ Row = 3; Col = 4; Letter = "D"; Value = "1.882.43" },
    @{ Row = 3; Col = 5; Letter = "E"; Value = "  +0.38%  " },
    @{ Row = 4; Col = 5; Letter = "E"; Value = "  +0.03%  " },
    @{ Row = 5; Col = 4; Letter = "D"; Value = "0.7133" },
    @{ Row = 5; Col = 5; Letter = "E"; Value = "  +0.14%  " },
    @{ Row = 6; Col = 4; Letter = "D"; Value = "242.86" },
    @{ Row = 6; Col = 5; Letter = "E"; Value = "  +0.38%  " },
    @{ Row = 8; Col = 4; Letter = "D"; Value = "0.08029" },
    @{ Row = 8; Col = 5; Letter = "E"; Value = "  +3.86%  " },
    @{ Row = 9; Col = 4; Letter = "D"; Value = "0.3128" },
    @{ Row = 9; Col = 5; Letter = "E"; Value = "  +0.91%  " },
    @{ Row = 10; Col = 4; Letter = "D"; Value = "25.19" },
    @{ Row = 10; Col = 5; Letter = "E"; Value = "  +1.43%  " },
    @{ Row = 11; Col = 4; Letter = "D"; Value = "0.08342" },
    @{ Row = 11; Col = 5; Letter = "E"; Value = "  -2.13%  " },
    @{ Row = 12; Col = 4; Letter = "D"; Value = "1.894.63" },
    @{ Row = 12; Col = 5; Letter = "E"; Value = "  +0.92%  " },
    @{ Row = 13; Col = 5; Letter = "E"; Value = "  +0.76%  " },
    @{ Row = 14; Col = 2; Letter = "B"; Value = "Litecoin" },
    @{ Row = 14; Col = 3; Letter = "C"; Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc" },
    @{ Row = 14; Col = 4; Letter = "D"; Value = "94.77" },
    @{ Row = 14; Col = 5; Letter = "E"; Value = "  +3.54%  " },
    @{ Row = 15; Col = 2; Letter = "B"; Value = "Polygon" },
    @{ Row = 15; Col = 3; Letter = "C"; Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic" },
    @{ Row = 15; Col = 4; Letter = "D"; Value = "0.7194" },
    @{ Row = 15; Col = 5; Letter = "E"; Value = "  +1.36%  " },
    @{ Row = 16; Col = 4; Letter = "D"; Value = "6.342" },
    @{ Row = 16; Col = 5; Letter = "E"; Value = "  +5.57%  " },
    @{ Row = 17; Col = 4; Letter = "D"; Value = "0.000008539" },
    @{ Row = 17; Col = 5; Letter = "E"; Value = "  +4.04%  " },
    @{ Row = 18; Col = 4; Letter = "D"; Value = "29.388.08" },
    @{ Row = 18; Col = 5; Letter = "E"; Value = "  +0.25%  " },
    @{ Row = 19; Col = 4; Letter = "D"; Value = "242.55" },
    @{ Row = 19; Col = 5; Letter = "E"; Value = "  +0.39%  " },
    @{ Row = 20; Col = 4; Letter = "D"; Value = "2.141.85" },
    @{ Row = 20; Col = 5; Letter = "E"; Value = "  +0.09%  " },
    @{ Row = 21; Col = 5; Letter = "E"; Value = "  +0.17%  " },
    @{ Row = 23; Col = 4; Letter = "D"; Value = "7.865" },
    @{ Row = 23; Col = 5; Letter = "E"; Value = "  +0.79%  " },
    @{ Row = 24; Col = 5; Letter = "E"; Value = "  -0.01%  " },
    @{ Row = 25; Col = 4; Letter = "D"; Value = "0.1587" },
    @{ Row = 25; Col = 5; Letter = "E"; Value = "  -0.83%  " },
    @{ Row = 26; Col = 4; Letter = "D"; Value = "163.30" },
    @{ Row = 26; Col = 5; Letter = "E"; Value = "  -0.06%  " },
    @{ Row = 27; Col = 4; Letter = "D"; Value = "9.083" },
    @{ Row = 27; Col = 5; Letter = "E"; Value = "  +0.58%  " },
    @{ Row = 28; Col = 5; Letter = "E"; Value = "  +0.93%  " },
    @{ Row = 29; Col = 5; Letter = "E"; Value = "  -0.16%  " },
    @{ Row = 31; Col = 4; Letter = "D"; Value = "4.336" },
    @{ Row = 31; Col = 5; Letter = "E"; Value = "  +0.49%  " },
    @{ Row = 32; Col = 5; Letter = "E"; Value = "  -6.65%  " },
    @{ Row = 33; Col = 4; Letter = "D"; Value = "0.05389" },
    @{ Row = 33; Col = 5; Letter = "E"; Value = "  +2.32%  " },
    @{ Row = 34; Col = 5; Letter = "E"; Value = "  +0.86%  " },
    @{ Row = 36; Col = 4; Letter = "D"; Value = "0.7516" },
    @{ Row = 36; Col = 5; Letter = "E"; Value = "  +0.94%  " },
    @{ Row = 37; Col = 4; Letter = "D"; Value = "2.699" },
    @{ Row = 37; Col = 5; Letter = "E"; Value = "  +0.52%  " },
    @{ Row = 38; Col = 4; Letter = "D"; Value = "0.01890" },
    @{ Row = 38; Col = 5; Letter = "E"; Value = "  +1.21%  " },
    @{ Row = 39; Col = 4; Letter = "D"; Value = "1.285.72" },
    @{ Row = 39; Col = 5; Letter = "E"; Value = "  +8.85%  " },
    @{ Row = 40; Col = 4; Letter = "D"; Value = "2.746" },
    @{ Row = 40; Col = 5; Letter = "E"; Value = "  +0.99%  " },
    @{ Row = 41; Col = 4; Letter = "D"; Value = "6.593" },
    @{ Row = 41; Col = 5; Letter = "E"; Value = "  +3.26%  " },
    @{ Row = 42; Col = 4; Letter = "D"; Value = "0.9162" },
    @{ Row = 42; Col = 5; Letter = "E"; Value = "  +3.39%  " },
    @{ Row = 43; Col = 4; Letter = "D"; Value = "74.65" },
    @{ Row = 44; Col = 4; Letter = "D"; Value = "111.96" },
    @{ Row = 44; Col = 5; Letter = "E"; Value = "  +5.27%  " },
    @{ Row = 45; Col = 4; Letter = "D"; Value = "1.001" },
    @{ Row = 45; Col = 5; Letter = "E"; Value = "  +0.03%  " },
    @{ Row = 46; Col = 5; Letter = "E"; Value = "  +6.77%  " },
    @{ Row = 47; Col = 4; Letter = "D"; Value = "2.031.68" },
    @{ Row = 47; Col = 5; Letter = "E"; Value = "  +0.06%  " },
    @{ Row = 48; Col = 4; Letter = "D"; Value = "1.810" },
    @{ Row = 48; Col = 5; Letter = "E"; Value = "  -0.06%  " },
    @{ Row = 49; Col = 4; Letter = "D"; Value = "0.5222" },
    @{ Row = 49; Col = 5; Letter = "E"; Value = "  +0.25%  " },
    @{ Row = 50; Col = 4; Letter = "D"; Value = "9.541" },
    @{ Row = 50; Col = 5; Letter = "E"; Value = "  +1.61%  " },
    @{ Row = 51; Col = 4; Letter = "D"; Value = "0.4393" },
    @{ Row = 51; Col = 5; Letter = "E"; Value = "  +1.77%  " },

)

foreach ($change in $changes) {
    Set-TextValue $change.Row $change.Col $change.Value
}
